$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Currency" column header
$ws.Range("E1").Value = "Currency"

# Add Currency value (AUD) before updating the Route value so that the
# shared-strings table picks up new unique strings in the same order as
# the target workbook (Currency, AUD, HIR-AKL).
$ws.Range("E2").Value = "AUD"

# Update Route value from HIR-BNE to HIR-AKL
$ws.Range("B2").Value = "HIR-AKL"

# Update selection to match target state
$ws.Range("B2").Select()
